$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'214"
$ws.Range("D2").Value = "'576228.00"
$ws.Range("C3").Value = "'1117"
$ws.Range("D3").Value = "'3876213.01"
$ws.Range("C4").Value = "'450"
$ws.Range("D4").Value = "'2119230.78"
$ws.Range("C6").Value = "'38"
$ws.Range("D6").Value = "'318643.82"
$ws.Range("C8").Value = "'52"
$ws.Range("D8").Value = "'113000.00"
$ws.Range("C15").Value = "'109"
$ws.Range("D15").Value = "'306126.38"
$ws.Range("C16").Value = "'447"
$ws.Range("D16").Value = "'1447153.23"
$ws.Range("C21").Value = "'58"
$ws.Range("D21").Value = "'143906.00"
$ws.Range("C22").Value = "'346"
$ws.Range("D22").Value = "'1164026.51"
$ws.Range("C23").Value = "'129"
$ws.Range("D23").Value = "'623660.00"
$ws.Range("C24").Value = "'47"
$ws.Range("D24").Value = "'310949.43"
$ws.Range("C35").Value = "'638"
$ws.Range("D35").Value = "'2279254.64"
$ws.Range("C36").Value = "'279"
$ws.Range("D36").Value = "'1625428.97"
$ws.Range("C38").Value = "'34"
$ws.Range("D38").Value = "'329055.00"
$ws.Range("C46").Value = "'41"
$ws.Range("D46").Value = "'149357.84"
$ws.Range("C47").Value = "'118"
$ws.Range("D47").Value = "'561404.84"
$ws.Range("C48").Value = "'65"
$ws.Range("D48").Value = "'442474.00"
$ws.Range("C51").Value = "'24"
$ws.Range("D51").Value = "'84174.00"
$ws.Range("C52").Value = "'121"
$ws.Range("D52").Value = "'385054.84"
$ws.Range("C53").Value = "'720"
$ws.Range("D53").Value = "'3067323.08"
$ws.Range("C54").Value = "'307"
$ws.Range("D54").Value = "'1583617.74"
$ws.Range("C55").Value = "'118"
$ws.Range("D55").Value = "'827097.18"
$ws.Range("C61").Value = "'748"
$ws.Range("D61").Value = "'4119097.06"
$ws.Range("C76").Value = "'14"
$ws.Range("D76").Value = "'70000.00"
$ws.Range("C78").Value = "'102"
$ws.Range("D78").Value = "'278639.87"
$ws.Range("C79").Value = "'436"
$ws.Range("D79").Value = "'1548208.84"
$ws.Range("C80").Value = "'168"
$ws.Range("D80").Value = "'826987.18"
$ws.Range("C81").Value = "'54"
$ws.Range("D81").Value = "'316844.67"
$ws.Range("C82").Value = "'17"
$ws.Range("D82").Value = "'160069.00"
$ws.Range("C83").Value = "'17"
$ws.Range("D83").Value = "'34000.00"
$ws.Range("C85").Value = "'956"
$ws.Range("D85").Value = "'3450600.16"
$ws.Range("C86").Value = "'364"
$ws.Range("D86").Value = "'1874553.44"
$ws.Range("C91").Value = "'291"
$ws.Range("D91").Value = "'867691.74"
$ws.Range("C92").Value = "'1103"
$ws.Range("D92").Value = "'3602501.20"
$ws.Range("C93").Value = "'401"
$ws.Range("D93").Value = "'1806017.61"
$ws.Range("C94").Value = "'123"
$ws.Range("D94").Value = "'608452.47"
$ws.Range("C95").Value = "'30"
$ws.Range("D95").Value = "'223000.00"
$ws.Range("C96").Value = "'61"
$ws.Range("D96").Value = "'122000.00"
$ws.Range("C105").Value = "'512"
$ws.Range("D105").Value = "'2285195.33"
$ws.Range("C106").Value = "'143"
$ws.Range("D106").Value = "'888996.00"
